$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1875
$ws.Range("C2").Value = 0.55
$ws.Range("J2").Value = 0.02142857142857143
$ws.Range("P2").Value = 0.1321428571428571
$ws.Range("S2").Value = 0.1089285714285714
$ws.Range("B3").Value = 0.0220820189274448
$ws.Range("C3").Value = 0.01261829652996845
$ws.Range("J3").Value = 0.05362776025236593
$ws.Range("P3").Value = 0.7192429022082019
$ws.Range("S3").Value = 0.1924290220820189
$ws.Range("J4").Value = 0.1016949152542373
$ws.Range("P4").Value = 0.7288135593220338
$ws.Range("S4").Value = 0.1694915254237288
$ws.Range("B6").Value = 0.05808080808080808
$ws.Range("D6").Value = 0.002525252525252525
$ws.Range("F6").Value = 0.08333333333333333
$ws.Range("J6").Value = 0.2853535353535354
$ws.Range("O6").Value = 0.01515151515151515
$ws.Range("Q6").Value = 0.1893939393939394
$ws.Range("R6").Value = 0.06060606060606061
$ws.Range("S6").Value = 0.3055555555555556
$ws.Range("B7").Value = 0.1165311653116531
$ws.Range("D7").Value = 0.01626016260162602
$ws.Range("F7").Value = 0.04607046070460705
$ws.Range("J7").Value = 0.1544715447154472
$ws.Range("O7").Value = 0.02710027100271003
$ws.Range("Q7").Value = 0.1680216802168022
$ws.Range("R7").Value = 0.08401084010840108
$ws.Range("S7").Value = 0.3875338753387534
$ws.Range("B8").Value = 0.12
$ws.Range("D8").Value = 0.01655172413793103
$ws.Range("F8").Value = 0.05931034482758621
$ws.Range("J8").Value = 0.1351724137931034
$ws.Range("O8").Value = 0.02206896551724138
$ws.Range("Q8").Value = 0.1572413793103448
$ws.Range("R8").Value = 0.08413793103448276
$ws.Range("S8").Value = 0.4055172413793103
$ws.Range("B9").Value = 0.1357615894039735
$ws.Range("D9").Value = 0.003311258278145695
$ws.Range("F9").Value = 0.05298013245033113
$ws.Range("J9").Value = 0.152317880794702
$ws.Range("O9").Value = 0.006622516556291391
$ws.Range("Q9").Value = 0.1754966887417219
$ws.Range("R9").Value = 0.08940397350993377
$ws.Range("S9").Value = 0.3841059602649007
$ws.Range("B10").Value = 0.1067615658362989
$ws.Range("D10").Value = 0.01912811387900356
$ws.Range("F10").Value = 0.06272241992882563
$ws.Range("J10").Value = 0.1534697508896797
$ws.Range("O10").Value = 0.02935943060498221
$ws.Range("Q10").Value = 0.202846975088968
$ws.Range("R10").Value = 0.0809608540925267
$ws.Range("S10").Value = 0.3447508896797153
$ws.Range("G11").Value = 0.1626297577854671
$ws.Range("J11").Value = 0.08304498269896193
$ws.Range("K11").Value = 0.2076124567474048
$ws.Range("L11").Value = 0.5380622837370242
$ws.Range("S11").Value = 0.00865051903114187
$ws.Range("G12").Value = 0.7292307692307692
$ws.Range("J12").Value = 0.2030769230769231
$ws.Range("K12").Value = 0.009230769230769232
$ws.Range("L12").Value = 0.03076923076923077
$ws.Range("S12").Value = 0.02769230769230769
$ws.Range("F13").Value = 0.01587301587301587
$ws.Range("G13").Value = 0.7301587301587301
$ws.Range("J13").Value = 0.2063492063492063
$ws.Range("S13").Value = 0.04761904761904762
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.01515151515151515
$ws.Range("H15").Value = 0.1616161616161616
$ws.Range("I15").Value = 0.06565656565656566
$ws.Range("J15").Value = 0.3383838383838384
$ws.Range("K15").Value = 0.08333333333333333
$ws.Range("M15").Value = 0.0101010101010101
$ws.Range("O15").Value = 0.06565656565656566
$ws.Range("S15").Value = 0.2601010101010101
$ws.Range("F16").Value = 0.02140672782874618
$ws.Range("H16").Value = 0.1712538226299694
$ws.Range("I16").Value = 0.06116207951070336
$ws.Range("J16").Value = 0.4587155963302753
$ws.Range("K16").Value = 0.1100917431192661
$ws.Range("M16").Value = 0.02140672782874618
$ws.Range("O16").Value = 0.04281345565749235
$ws.Range("S16").Value = 0.1131498470948012
$ws.Range("F17").Value = 0.02127659574468085
$ws.Range("H17").Value = 0.1555851063829787
$ws.Range("I17").Value = 0.0851063829787234
$ws.Range("J17").Value = 0.4148936170212766
$ws.Range("K17").Value = 0.1170212765957447
$ws.Range("M17").Value = 0.0199468085106383
$ws.Range("O17").Value = 0.05851063829787234
$ws.Range("S17").Value = 0.1276595744680851
$ws.Range("F18").Value = 0.02476780185758514
$ws.Range("H18").Value = 0.1764705882352941
$ws.Range("I18").Value = 0.09907120743034056
$ws.Range("J18").Value = 0.4334365325077399
$ws.Range("K18").Value = 0.1052631578947368
$ws.Range("M18").Value = 0.01238390092879257
$ws.Range("O18").Value = 0.04334365325077399
$ws.Range("S18").Value = 0.1052631578947368
$ws.Range("F19").Value = 0.03111327175498298
$ws.Range("H19").Value = 0.2134175984443364
$ws.Range("I19").Value = 0.07632474477394263
$ws.Range("J19").Value = 0.3534273213417599
$ws.Range("K19").Value = 0.1239669421487603
$ws.Range("M19").Value = 0.0160427807486631
$ws.Range("N19").Value = 0.0009722897423432182
$ws.Range("O19").Value = 0.07681088964511425
$ws.Range("S19").Value = 0.1079241614000972
